$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new, empty paragraph right after the existing blank paragraph
#    that follows "Vomits blood" (i.e. immediately before the "Self-care
#    strategies" paragraph).
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Self-care strategies*") {
        $prevPara = $p.Previous()
        $insertRange = $prevPara.Range.Duplicate()
        $insertRange.InsertParagraphAfter()
        $newPara = $prevPara.Next()
        $newPara.Range.Text = ""
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Update the self-care copy: "water" -> "fluids", and drop the trailing
#    period.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Give your child plenty of water*") {
        $p.Range.Text = "Give your child plenty of fluids if the abdominal pain is accompanied by diarrhea or constipation"
        break
    }
}
